# "Generate Report for Handback" — fills in the Latest Target File / Latest
# Handback File / Latest Handback DateTime columns for both locale sheets
# (zh-cn, de-de), flips the Status text from "Ready for handoff" to
# "Handed back: in sync with en-US", adds the matching hyperlinks for the
# newly-populated "Latest Target File" cells, and widens a few columns that
# now hold longer text.

$wb = $excel.ActiveWorkbook

$baseUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/8191ff7835398d93ec5a2512e0a5ca224db7b29c/e2e/"
$file1 = "3abb78c5-e03d-42e6-a5a1-4531fd5de059.md"
$file2 = "54f14ea3-ff94-4d72-8bbf-ea4d62077458.md"
$url1 = $baseUrl + $file1
$url2 = $baseUrl + $file2

$zhXlf1 = "3abb78c5-e03d-42e6-a5a1-4531fd5de059.77949e912d2d57e882c0663f93a2b2c2d44851ce.zh-cn.xlf"
$zhXlf2 = "54f14ea3-ff94-4d72-8bbf-ea4d62077458.f705c999e088170bc7e7bd5deb9b62e2c7cd06f1.zh-cn.xlf"
$deXlf1 = "3abb78c5-e03d-42e6-a5a1-4531fd5de059.77949e912d2d57e882c0663f93a2b2c2d44851ce.de-de.xlf"
$deXlf2 = "54f14ea3-ff94-4d72-8bbf-ea4d62077458.f705c999e088170bc7e7bd5deb9b62e2c7cd06f1.de-de.xlf"

$zhHandbackDate = "2016-08-13 05:15:44"
$deHandbackDate = "2016-08-13 05:15:55"

$newStatus = "Handed back: in sync with en-US"

# ── Overview sheet: status column (E/F) text now reads longer, so the
#    columns need to be widened to keep the same visual fit. ────────────────
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus
$overview.Columns.Item(5).ColumnWidth = 29.166666666666668
$overview.Columns.Item(6).ColumnWidth = 29.166666666666668

# ── zh-cn sheet ──────────────────────────────────────────────────────────
$zh = $wb.Worksheets.Item("zh-cn")

$zh.Range("C2").Value = $newStatus
$zh.Range("C3").Value = $newStatus

$zh.Range("I2").Value = $file1
$zh.Range("I2").Style = "HyperLink"
$zh.Range("J2").Value = $zhXlf1
$zh.Range("K2").Value = $zhHandbackDate

$zh.Range("I3").Value = $file2
$zh.Range("I3").Style = "HyperLink"
$zh.Range("J3").Value = $zhXlf2
$zh.Range("K3").Value = $zhHandbackDate

$zh.Hyperlinks.Delete()
$zh.Hyperlinks.Add($zh.Range("A2"), $url1, "", "", $file1)
$zh.Hyperlinks.Add($zh.Range("I2"), $url1, "", "", $file1)
$zh.Hyperlinks.Add($zh.Range("A3"), $url2, "", "", $file2)
$zh.Hyperlinks.Add($zh.Range("I3"), $url2, "", "", $file2)

$zh.Columns.Item(3).ColumnWidth = 29.166666666666668
$zh.Columns.Item(9).ColumnWidth = 39.166666666666664
$zh.Columns.Item(10).ColumnWidth = 39.166666666666664

# ── de-de sheet ──────────────────────────────────────────────────────────
$de = $wb.Worksheets.Item("de-de")

$de.Range("C2").Value = $newStatus
$de.Range("C3").Value = $newStatus

$de.Range("I2").Value = $file1
$de.Range("I2").Style = "HyperLink"
$de.Range("J2").Value = $deXlf1
$de.Range("K2").Value = $deHandbackDate

$de.Range("I3").Value = $file2
$de.Range("I3").Style = "HyperLink"
$de.Range("J3").Value = $deXlf2
$de.Range("K3").Value = $deHandbackDate

$de.Hyperlinks.Delete()
$de.Hyperlinks.Add($de.Range("A2"), $url1, "", "", $file1)
$de.Hyperlinks.Add($de.Range("I2"), $url1, "", "", $file1)
$de.Hyperlinks.Add($de.Range("A3"), $url2, "", "", $file2)
$de.Hyperlinks.Add($de.Range("I3"), $url2, "", "", $file2)

$de.Columns.Item(3).ColumnWidth = 29.166666666666668
$de.Columns.Item(9).ColumnWidth = 39.166666666666664
$de.Columns.Item(10).ColumnWidth = 39.166666666666664

Write-Host "Handback report generated."
